$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2, 250768, 453260),
    @(3, 1003593, 1782688),
    @(4, 2255654, 4077636),
    @(5, 4002551, 7107394),
    @(6, 6246158, 11094581),
    @(7, 8998411, 15763670),
    @(8, 12269930, 21612353),
    @(9, 16023725, 28111146),
    @(10, 20260468, 35240291),
    @(11, 25018073, 43538819),
    @(12, 30222063, 53007276),
    @(13, 36012208, 62442496),
    @(14, 42239015, 73601961),
    @(15, 48973565, 86172707),
    @(16, 56208277, 98356119),
    @(17, 64043158, 111128785),
    @(18, 72246481, 126022208),
    @(19, 80999775, 141312817),
    @(20, 90225541, 156804226),
    @(21, 99984571, 174203609),
    @(22, 110289460, 192302282),
    @(23, 120956741, 210411953),
    @(24, 132101619, 227473285),
    @(25, 144091178, 251880984),
    @(26, 156291361, 273094577),
    @(27, 169198551, 294246289),
    @(28, 182223143, 318103043),
    @(29, 196098204, 341762631),
    @(30, 210417647, 366048123),
    @(31, 225291240, 391829383),
    @(32, 240299213, 419469505),
    @(33, 256151588, 446190977),
    @(34, 272186802, 474372755),
    @(35, 288884207, 502494080),
    @(36, 306703506, 532409922),
    @(37, 323968637, 560984995),
    @(38, 341890893, 590764442),
    @(39, 361347963, 628461309),
    @(40, 380397598, 661364272),
    @(41, 399875610, 694500135),
    @(42, 420355652, 730470567),
    @(43, 441090503, 765262131),
    @(44, 462668697, 803366141),
    @(45, 484128793, 830996811),
    @(46, 506080649, 863139613),
    @(47, 528931551, 902397807),
    @(48, 552936897, 942953781),
    @(49, 575691892, 982592004),
    @(50, 600338646, 1024352685),
    @(51, 624902744, 1066160511),
    @(52, 650719630, 1109685106),
    @(53, 676422807, 1153831765),
    @(54, 702176007, 1198926462),
    @(55, 728990471, 1245098520),
    @(56, 756900881, 1290884783),
    @(57, 784706785, 1339577889),
    @(58, 812312953, 1386343785),
    @(59, 840754819, 1435025148),
    @(60, 870438528, 1484677034),
    @(61, 899846055, 1535214709),
    @(62, 930153440, 1587653969),
    @(63, 961206936, 1639596509),
    @(64, 992000301, 1692447016),
    @(65, 1023805767, 1745980417),
    @(66, 1056496482, 1800771297),
    @(67, 1088085636, 1855797557),
    @(68, 1122047444, 1913924446),
    @(69, 1156811322, 1972220061),
    @(70, 1190166004, 2029359055),
    @(71, 1224092269, 2086829539),
    @(72, 1260588338, 2149250961),
    @(73, 1295866967, 2209787418),
    @(74, 1331713664, 2273871443),
    @(75, 1369252716, 2333709003),
    @(76, 1406129277, 2398702629),
    @(77, 1444126793, 2463910201),
    @(78, 1483152271, 2529369010),
    @(79, 1520780002, 2593963320),
    @(80, 1559792607, 2660482691),
    @(81, 1599423580, 2726735362),
    @(82, 1640193311, 2797815129),
    @(83, 1680034237, 2865976252),
    @(84, 1723123944, 2937501925),
    @(85, 1764294091, 3008133632),
    @(86, 1806096468, 3078032103),
    @(87, 1849373479, 3152725145),
    @(88, 1891844672, 3224970352),
    @(89, 1935630918, 3300921352),
    @(90, 1980006919, 3377607530),
    @(91, 2024927528, 3451586382),
    @(92, 2070133780, 3528114294),
    @(93, 2115103978, 3605654255),
    @(94, 2161945660, 3687141563),
    @(95, 2208588182, 3764647295),
    @(96, 2255022149, 3844976955),
    @(97, 2304058704, 3926944005),
    @(98, 2352987386, 4011530009),
    @(99, 2400495939, 4090189189),
    @(100, 2449798006, 4172962107),
    @(101, 2499446093, 4260845392),
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
}